$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Kitchen & Dining / Kitchen & Dining
$ws.Range("C2").Value = 8264

# Row 3: Home Decor / Home Decor
$ws.Range("C3").Value = 8252

# Row 4: Stationery & Office / Stationery & Office
$ws.Range("C4").Value = 3630

# Row 5: Home Decor / Kitchen & Dining
$ws.Range("A5").Value = "Home Decor"
$ws.Range("B5").Value = "Kitchen & Dining"
$ws.Range("C5").Value = 3544

# Row 6: Kitchen & Dining / Home Decor
$ws.Range("A6").Value = "Kitchen & Dining"
$ws.Range("B6").Value = "Home Decor"
$ws.Range("C6").Value = 3544

# Row 7: Seasonal & Holidays / Seasonal & Holidays
$ws.Range("A7").Value = "Seasonal & Holidays"
$ws.Range("B7").Value = "Seasonal & Holidays"
$ws.Range("C7").Value = 3072

# Row 8: Kids & Toys / Kids & Toys
$ws.Range("A8").Value = "Kids & Toys"
$ws.Range("B8").Value = "Kids & Toys"
$ws.Range("C8").Value = 1908

# Row 9: Home Decor / Seasonal & Holidays
$ws.Range("A9").Value = "Home Decor"
$ws.Range("B9").Value = "Seasonal & Holidays"
$ws.Range("C9").Value = 1856

# Row 10: Seasonal & Holidays / Home Decor
$ws.Range("A10").Value = "Seasonal & Holidays"
$ws.Range("B10").Value = "Home Decor"
$ws.Range("C10").Value = 1856

# Row 11: Home Decor / Stationery & Office
$ws.Range("A11").Value = "Home Decor"
$ws.Range("B11").Value = "Stationery & Office"
$ws.Range("C11").Value = 1843
